$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(101, "Yes", 21, "11/03/2022 - Board0 - NO, RH"),
    @(102, "Yes", 22, "11/03/2022 - Board1 - NO, RH"),
    @(103, "Yes", 21, "11/04/2022 - Board0 - NO, RH"),
    @(104, "Yes", 21, "11/07/2022 - Board0 - N2O, RH"),
    @(105, "Yes", 21, "11/07/2022 - Board0 - N2O, RH Permselect"),
    @(106, "Yes", 21, "11/07/2022 - Board0 - N2O, RH Permselect, PostSCUID"),
    @(107, "Yes", 22, "11/08/2022 - Board0 - N2O, RH Permselect, Air"),
    @(108, "Yes", 22, "11/08/2022 - Board0 - N2O, RH Permselect")
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Style = "Normal"
}
